$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the email values (column I and J hold the same email) ---
$ws.Range("I2").Value = "plhm@plaka.com"
$ws.Range("J2").Value = "plhm@plaka.com"

$ws.Range("I3").Value = "okmh@okaaj.com"
$ws.Range("J3").Value = "okmh@okaaj.com"

$ws.Range("I4").Value = "ijhh@waska.com"
$ws.Range("J4").Value = "ijhh@waska.com"

$ws.Range("I5").Value = "oosh@gagap.com"
$ws.Range("J5").Value = "oosh@gagap.com"

$ws.Range("I6").Value = "asvh@waear.com"
$ws.Range("J6").Value = "asvh@waear.com"

$ws.Range("I7").Value = "qplh@sanax.com"
$ws.Range("J7").Value = "qplh@sanax.com"

# --- Remove the wrap-text formatting on column D (address) ---
# Copy the plain text-format style (no wrap) from a neighbouring cell
# so the previously-used "wrap text" style becomes unused, instead of
# creating a brand-new style.
$null = $ws.Range("E2").Copy()
$null = $ws.Range("D2:D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows no longer need the extra height that wrapping required ---
$null = $ws.Rows("2:7").AutoFit()

# --- Column D should now size itself to its (no-longer-wrapped) contents ---
$null = $ws.Columns("D").AutoFit()

# --- Page orientation set to portrait ---
$ws.PageSetup.Orientation = 1

# --- Selection moves to the full data range ---
$null = $ws.Range("A2:J7").Select()
